$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "theta_threshold_range" row (row 5) entirely, shifting rows below it up.
$ws.Rows.Item(5).Delete()

# Update values on the remaining rows to match the new data set.
$ws.Range("B2").Value = 5.7
$ws.Range("C2").Value = 10.6

$ws.Range("B3").Value = 5.6
$ws.Range("C3").Value = 9.3000000000000007

$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 1.3

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Update the selected cell / window to match the saved view state.
$ws.Range("B4").Select()
